# Added BIS Billing V3 CII binding
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Document Type")
$ws.Activate()

# Copy the formatting of the row above (same style pattern PEPPOL BIS Billing V3
# rows already use: A/B = style 4, C/D = style 6) onto the new row 35.
$ws.Range("A31:D31").Copy()
$ws.Range("A35:D35").PasteSpecial(-4122) # xlPasteFormats

# Row 35: new PEPPOL BIS Billing V3 CII document identifier binding
$ws.Cells.Item(35, 1).Value = "PEPPOL BIS Billing V3"
$ws.Cells.Item(35, 2).Value = "urn:un:unece:uncefact:data:standard:CrossIndustryInvoice:100::CrossIndustryInvoice##urn:cen.eu:en16931:2017#compliant#urn:fdc:peppol.eu:2017:poacc:billing:3.0::D16B"
$ws.Cells.Item(35, 3).Value = 3
$ws.Cells.Item(35, 4).Formula = "=FALSE"

$ws.Range("A35:D35").RowHeight = 30

# Update scroll position / active selection to reflect the appended row
$excel.Goto($ws.Range("A23"), $true)
$ws.Range("B36").Select()
